$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 8) continuing the data table: Sr No = 4, Detail = "This is Version 4 of File"
$ws.Range("B8").Value = 4
$ws.Range("C8").Value = "This is Version 4 of File"

# Update the active selection to match the post-edit state (C9, just below the new data)
$ws.Range("C9").Select()
